$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.691.15"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.252.67"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.648"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.444"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0969"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "2.595.78"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.837"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "2.263.90"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "43.660.53"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "0.0₃0972"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +31.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.136"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0680"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0253"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0948"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.453.87"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("B49").Value = "TerraClassic"
$ws.Range("C49").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000206"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("B51").Value = "Celestia"
$ws.Range("C51").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.41%  "
